# Migration to Automation-Org/TestCases-maintenance/WIP-RMA TestCases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column H: data rows with "a811K0000004fpN" and header "Labor Booking User"
# (bold header style, matching A1:F1). Data cells written first so the shared-string
# table order matches (a811K0000004fpN before Labor Booking User).
$ws.Range("H2").Value = "a811K0000004fpN"
$ws.Range("H3").Value = "a811K0000004fpN"

$ws.Range("H1").Value = "Labor Booking User"
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Update the selection to match the post-edit state (active cell H1, single cell selected)
$ws.Range("H1").Select()
